$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 99.875
$ws.Range("I9").Value = 61.761906
$ws.Range("K9").Value = 61.761906
$ws.Range("M9").Value = 107.238094
$ws.Range("H69").Value = 3378.6047
$ws.Range("J69").Value = 3300
$ws.Range("L69").Value = 9900
$ws.Range("N69").Value = -11648
$ws.Range("H72").Value = 3378.6047
$ws.Range("J72").Value = 3300
$ws.Range("L72").Value = 29700
$ws.Range("N72").Value = -38436
$ws.Range("H75").Value = 23328.5
$ws.Range("J75").Value = 23328.5
$ws.Range("L75").Value = 23328.5
$ws.Range("N75").Value = -25200.5
$ws.Range("H78").Value = 23328.5
$ws.Range("J78").Value = 23328.5
$ws.Range("L78").Value = 69985.5
$ws.Range("N78").Value = -79345.5
$ws.Range("H106").Value = 29471568
$ws.Range("I106").Value = 67710.60000000001
$ws.Range("J106").Value = 250000500
$ws.Range("K106").Value = 67710.60000000001
$ws.Range("L106").Value = 250000500
$ws.Range("M106").Value = -67079.60000000001
$ws.Range("N106").Value = -250001762
$ws.Range("H116").Value = 1593.0476
$ws.Range("I116").Value = 1377.091
$ws.Range("K116").Value = 1377.091
$ws.Range("M116").Value = 2064.909
$ws.Range("H132").Value = 5496312
$ws.Range("I132").Value = 6213109
$ws.Range("J132").Value = 868.6667
$ws.Range("K132").Value = 18639327
$ws.Range("L132").Value = 2606.0001
$ws.Range("M132").Value = -18636797
$ws.Range("N132").Value = -7666.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2945.1084
$ws.Range("I32").Value = 2803.2534
$ws.Range("J32").Value = 4275
$ws.Range("K32").Value = 2803.2534
$ws.Range("L32").Value = 4275
$ws.Range("M32").Value = -2516.2534
$ws.Range("N32").Value = -4849
$ws.Range("H76").Value = 25145
$ws.Range("J76").Value = 25145
$ws.Range("L76").Value = 25145
$ws.Range("N76").Value = -25821
$ws.Range("H79").Value = 25145
$ws.Range("J79").Value = 25145
$ws.Range("L79").Value = 25145
$ws.Range("N79").Value = -27485

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372
$ws.Range("H64").Value = 375.8889
$ws.Range("J64").Value = 359.57144
$ws.Range("L64").Value = 359.57144
$ws.Range("N64").Value = -809.5714399999999
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864
$ws.Range("H67").Value = 375.8889
$ws.Range("J67").Value = 359.57144
$ws.Range("L67").Value = 359.57144
$ws.Range("N67").Value = -1919.57144
$ws.Range("H80").Value = 707.6923
$ws.Range("J80").Value = 742.4286
$ws.Range("L80").Value = 742.4286
$ws.Range("N80").Value = -2738.4286
$ws.Range("H83").Value = 707.6923
$ws.Range("J83").Value = 742.4286
$ws.Range("L83").Value = 3712.143
$ws.Range("N83").Value = -13696.143

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3970522.8
$ws.Range("I31").Value = 1749.6086
$ws.Range("J31").Value = 8774827
$ws.Range("K31").Value = 1749.6086
$ws.Range("L31").Value = 8774827
$ws.Range("M31").Value = -1454.6086
$ws.Range("N31").Value = -8775417
$ws.Range("H34").Value = 3970522.8
$ws.Range("I34").Value = 1749.6086
$ws.Range("J34").Value = 8774827
$ws.Range("K34").Value = 1749.6086
$ws.Range("L34").Value = 8774827
$ws.Range("M34").Value = -1547.6086
$ws.Range("N34").Value = -8775231
$ws.Range("H62").Value = 3075
$ws.Range("J62").Value = 3266.6667
$ws.Range("L62").Value = 3266.6667
$ws.Range("N62").Value = -4514.6667
$ws.Range("H65").Value = 3075
$ws.Range("J65").Value = 3266.6667
$ws.Range("L65").Value = 16333.3335
$ws.Range("N65").Value = -22573.3335
$ws.Range("H94").Value = 344.4
$ws.Range("I94").Value = 305.5
$ws.Range("J94").Value = 500
$ws.Range("K94").Value = 305.5
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = 145.5
$ws.Range("N94").Value = -1402
$ws.Range("H132").Value = 2165.2354
$ws.Range("I132").Value = 1343.7142
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 4031.1426
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -1501.1426
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 1166.96
$ws.Range("I134").Value = 1165.3334
$ws.Range("J134").Value = 1171.1428
$ws.Range("K134").Value = 3496.0002
$ws.Range("L134").Value = 3513.4284
$ws.Range("M134").Value = -961.0001999999999
$ws.Range("N134").Value = -8583.428400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 881.25
$ws.Range("I69").Value = 675
$ws.Range("J69").Value = 1500
$ws.Range("K69").Value = 2025
$ws.Range("L69").Value = 4500
$ws.Range("M69").Value = -1214
$ws.Range("N69").Value = -6122
$ws.Range("H72").Value = 881.25
$ws.Range("I72").Value = 675
$ws.Range("J72").Value = 1500
$ws.Range("K72").Value = 6075
$ws.Range("L72").Value = 13500
$ws.Range("M72").Value = -2019
$ws.Range("N72").Value = -21612
$ws.Range("H131").Value = 3763.7026
$ws.Range("J131").Value = 920.62964
$ws.Range("L131").Value = 2761.88892
$ws.Range("N131").Value = -12841.88892

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1968.0625
$ws.Range("I126").Value = 1792.6
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 5377.799999999999
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -2907.799999999999
$ws.Range("N126").Value = -18740
$ws.Range("H132").Value = 2655
$ws.Range("I132").Value = 2345.276
$ws.Range("K132").Value = 7035.828
$ws.Range("M132").Value = -4505.828
$ws.Range("H141").Value = 39838.168
$ws.Range("J141").Value = 39838.168
$ws.Range("L141").Value = 39838.168
$ws.Range("N141").Value = -50198.168

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 15716.667
$ws.Range("I64").Value = 15000
$ws.Range("J64").Value = 16075
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = 16075
$ws.Range("M64").Value = -14775
$ws.Range("N64").Value = -16525
$ws.Range("H67").Value = 15716.667
$ws.Range("I67").Value = 15000
$ws.Range("J67").Value = 16075
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 16075
$ws.Range("M67").Value = -14220
$ws.Range("N67").Value = -17635

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184
$ws.Range("H93").Value = 32194.5
$ws.Range("J93").Value = 32194.5
$ws.Range("L93").Value = 32194.5
$ws.Range("N93").Value = -37186.5

Write-Output "Applied 179 cell updates across 8 sheets"